# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# symbol list with the latest scrape values.
#
# Values are formatted as free-text strings in the source workbook
# (percent signs, fixed decimal places, etc. are literal characters,
# not a numeric percentage format), so each target cell is forced to
# Text ("@") number format before the new literal is written. This
# keeps e.g. "0.1000" / "0.00000000750" from losing trailing zeros and
# keeps "-2.95%" stored as the literal string rather than being parsed
# into a numeric percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
}

# Row 2
Set-TextValue "D2" "293.68"
Set-TextValue "E2" "-2.95%"

# Row 3
Set-TextValue "D3" "31.08"
Set-TextValue "E3" "-2.74%"

# Row 4
Set-TextValue "D4" "4.891"
Set-TextValue "E4" "-2.17%"

# Row 5
Set-TextValue "D5" "0.07342"
Set-TextValue "E5" "-6.89%"

# Row 6
Set-TextValue "D6" "1.833"
Set-TextValue "E6" "-13.55%"

# Row 7
Set-TextValue "D7" "7.674"
Set-TextValue "E7" "-2.04%"

# Row 8
Set-TextValue "D8" "3.763"
Set-TextValue "E8" "-0.86%"

# Row 9
Set-TextValue "D9" "0.9059"
Set-TextValue "E9" "-2.46%"

# Row 10
Set-TextValue "D10" "0.1655"
Set-TextValue "E10" "-5.67%"

# Row 11
Set-TextValue "D11" "0.07541"
Set-TextValue "E11" "-5.10%"

# Row 12
Set-TextValue "D12" "0.08127"
Set-TextValue "E12" "-7.73%"

# Row 13
Set-TextValue "D13" "0.02988"
Set-TextValue "E13" "-4.47%"

# Row 14
Set-TextValue "D14" "0.1000"
Set-TextValue "E14" "-0.32%"

# Row 15
Set-TextValue "D15" "0.001495"
Set-TextValue "E15" "-0.97%"

# Row 16
Set-TextValue "D16" "0.005702"
Set-TextValue "E16" "-4.28%"

# Row 17
Set-TextValue "D17" "3.459"

# Row 18
Set-TextValue "E18" "-7.99%"

# Row 19
Set-TextValue "D19" "0.3282"
Set-TextValue "E19" "-0.36%"

# Row 20
Set-TextValue "E20" "1.26%"

# Row 21
Set-TextValue "D21" "4.348"
Set-TextValue "E21" "4.65%"

# Row 22
Set-TextValue "E22" "11.90%"

# Row 23
Set-TextValue "D23" "0.04485"

# Row 24
Set-TextValue "D24" "0.001226"
Set-TextValue "E24" "-0.75%"

# Row 25
Set-TextValue "D25" "0.004045"
Set-TextValue "E25" "-10.28%"

# Row 26
Set-TextValue "D26" "0.0001251"
Set-TextValue "E26" "0.12%"

# Row 39
Set-TextValue "D39" "0.01653"
Set-TextValue "E39" "-4.82%"

# Row 40
Set-TextValue "D40" "0.04391"
Set-TextValue "E40" "-7.97%"

# Row 41
Set-TextValue "D41" "0.007412"
Set-TextValue "E41" "1.10%"

# Row 42
Set-TextValue "D42" "0.1321"
Set-TextValue "E42" "-3.34%"

# Row 43
Set-TextValue "D43" "0.002091"
Set-TextValue "E43" "-10.57%"

# Row 44
Set-TextValue "D44" "0.01108"
Set-TextValue "E44" "1.34%"

# Row 45
Set-TextValue "D45" "0.00006010"
Set-TextValue "E45" "-0.97%"

# Row 46
Set-TextValue "D46" "0.00000000750"
Set-TextValue "E46" "0.11%"

# Row 47
Set-TextValue "D47" "2.149"
Set-TextValue "E47" "161.93%"

# Row 48
Set-TextValue "D48" "0.002401"
Set-TextValue "E48" "-29.31%"

# Row 49
Set-TextValue "D49" "0.00002101"
Set-TextValue "E49" "0.11%"

# Row 50
Set-TextValue "D50" "0.0002001"
Set-TextValue "E50" "0.11%"
